# "Fix level 4 and complete level 6"
#
# The underlying grid on the "Level8" sheet stores a level-design tile map
# as plain numbers. A batch of cells that were left at the placeholder
# value `1` are corrected to their real tile ids (17, 15, 14, 16, or 22
# for most of them). We also restore the author's last on-screen
# selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Tile-id corrections ----------------------------------------------
$fixes = @{
    "Q3"   = 22
    "Y3"   = 22
    "AG3"  = 22
    "AZ5"  = 22
    "AV8"  = 17
    "BI8"  = 15
    "BI9"  = 22
    "E13"  = 22
    "AQ13" = 22
    "U14"  = 22
    "AX15" = 22
    "BB17" = 22
    "X19"  = 22
    "BF19" = 22
    "E21"  = 22
    "J21"  = 22
    "AL21" = 22
    "AQ21" = 22
    "BJ21" = 22
    "U22"  = 22
    "AA22" = 22
    "X25"  = 22
    "BD25" = 22
    "N26"  = 22
    "E29"  = 22
    "AQ29" = 22
    "AU29" = 22
    "BE29" = 16
    "BA32" = 22
    "BF32" = 22
    "BA33" = 14
    "BD33" = 17
    "X34"  = 22
    "Y35"  = 22
    "AV36" = 17
    "BC38" = 22
    "Q41"  = 22
    "Y41"  = 22
    "AG41" = 22
}

foreach ($addr in $fixes.Keys) {
    $ws.Range($addr).Value = $fixes[$addr]
}

# --- Restore the saved view/selection ----------------------------------
$win = $excel.ActiveWindow
$ws.Range("AP31").Select()
$win.ScrollRow = 7
$win.ScrollColumn = 1
